# Apply updated odds values to Sheet1, matching the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.08
$ws.Range("R5").Value = 1.73

# Row 8
$ws.Range("I8").Value = 6.5
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("AG8").Value = 301
$ws.Range("AM8").Value = 51
$ws.Range("AS8").Value = 151
$ws.Range("AX8").Value = 34
$ws.Range("AZ8").Value = 126

# Row 10
$ws.Range("G10").Value = 1.7
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 5.5
$ws.Range("J10").Value = 2.38
$ws.Range("K10").Value = 2.1
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("U10").Value = 2.1
$ws.Range("V10").Value = 1.67
$ws.Range("X10").Value = 7
$ws.Range("AN10").Value = 3.5
$ws.Range("AO10").Value = 9

# Row 11
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
